$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '56.500.71'
Set-TextValue 'E2' '  -1.28%  '
Set-TextValue 'D3' '2.315.77'
Set-TextValue 'E3' '  -0.42%  '
Set-TextValue 'E4' '  -0.08%  '
Set-TextValue 'D5' '512.59'
Set-TextValue 'E5' '  -1.82%  '
Set-TextValue 'D6' '131.27'
Set-TextValue 'E6' '  -2.92%  '
Set-TextValue 'E7' '  +0.21%  '
Set-TextValue 'D8' '0.533'
Set-TextValue 'E8' '  -1.02%  '
Set-TextValue 'D9' '0.100'
Set-TextValue 'E9' '  -3.69%  '
Set-TextValue 'E10' '  -0.12%  '
Set-TextValue 'D11' '5.24'
Set-TextValue 'E11' '  -1.14%  '
Set-TextValue 'E12' '  -2.26%  '
Set-TextValue 'B13' 'Avalanche'
Set-TextValue 'C13' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D13' '23.47'
Set-TextValue 'E13' '  -2.15%  '
Set-TextValue 'B14' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D14' '2.728.34'
Set-TextValue 'E14' '  -0.90%  '
Set-TextValue 'D15' '56.461.33'
Set-TextValue 'E15' '  -0.83%  '
Set-TextValue 'E16' '  -1.98%  '
Set-TextValue 'D17' '2.325.11'
Set-TextValue 'E17' '  -0.81%  '
Set-TextValue 'D18' '10.33'
Set-TextValue 'E18' '  -1.65%  '
Set-TextValue 'D19' '327.54'
Set-TextValue 'E19' '  +1.35%  '
Set-TextValue 'D20' '4.13'
Set-TextValue 'E20' '  -2.25%  '
Set-TextValue 'D21' '6.72'
Set-TextValue 'E21' '  +1.52%  '
Set-TextValue 'D22' '0.999'
Set-TextValue 'E22' '  +0.02%  '
Set-TextValue 'D23' '61.22'
Set-TextValue 'E23' '  +0.76%  '
Set-TextValue 'D24' '0.164'
Set-TextValue 'E24' '  -0.90%  '
Set-TextValue 'D25' '8.58'
Set-TextValue 'E25' '  +7.55%  '
Set-TextValue 'E26' '  +0.67%  '
Set-TextValue 'D27' '1.31'
Set-TextValue 'E27' '  +0.45%  '
Set-TextValue 'D28' '167.52'
Set-TextValue 'E28' '  -0.55%  '
Set-TextValue 'D29' '1.67'
Set-TextValue 'E29' '  -3.67%  '
Set-TextValue 'D30' '0.0₃0717'
Set-TextValue 'E30' '  -4.00%  '
Set-TextValue 'D31' '6.08'
Set-TextValue 'E31' '  -1.96%  '
Set-TextValue 'D32' '18.28'
Set-TextValue 'E32' '  -0.46%  '
Set-TextValue 'E33' '  -0.02%  '
Set-TextValue 'E34' '  +0.11%  '
Set-TextValue 'E35' '  -1.59%  '
Set-TextValue 'D36' '3.91'
Set-TextValue 'E36' '  -3.10%  '
Set-TextValue 'D37' '0.884'
Set-TextValue 'E37' '  -4.82%  '
Set-TextValue 'B38' 'Stacks'
Set-TextValue 'C38' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D38' '1.56'
Set-TextValue 'E38' '  +0.11%  '
Set-TextValue 'B39' 'OKB'
Set-TextValue 'C39' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D39' '38.54'
Set-TextValue 'E39' '  +1.62%  '
Set-TextValue 'D40' '149.46'
Set-TextValue 'E40' '  +7.92%  '
Set-TextValue 'D41' '0.373'
Set-TextValue 'E41' '  -1.76%  '
Set-TextValue 'D42' '3.56'
Set-TextValue 'E42' '  -1.00%  '
Set-TextValue 'D43' '275.72'
Set-TextValue 'E43' '  -0.82%  '
Set-TextValue 'D44' '4.99'
Set-TextValue 'E44' '  -4.54%  '
Set-TextValue 'D45' '0.0925'
Set-TextValue 'E45' '  -0.92%  '
Set-TextValue 'D46' '0.0494'
Set-TextValue 'E46' '  -2.66%  '
Set-TextValue 'D47' '0.552'
Set-TextValue 'E47' '  -2.19%  '
Set-TextValue 'D48' '18.19'
Set-TextValue 'E48' '  +1.96%  '
Set-TextValue 'E49' '  -1.87%  '
Set-TextValue 'B50' 'Polygon'
Set-TextValue 'C50' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D50' '0.376'
Set-TextValue 'E50' '  -0.66%  '
Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '17.04'
Set-TextValue 'E51' '  +0.82%  '
